# Updates cryptos list values (price / 1h volume %, and a coin re-rank)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text in the sheet (values like "48.065.05"
# use "." as a thousands separator, and some end in a significant trailing
# zero, e.g. "20.20"). Prefix new D-column values with a leading apostrophe
# so Excel stores them as text instead of re-parsing them as numbers (which
# would silently drop the trailing zero / misread the separators).
$txt = "'"

# Row 2: Bitcoin
$ws.Range("D2").Value = $txt + '47.982.13'
$ws.Range("E2").Value = '  -0.09%  '

# Row 3: Ethereum
$ws.Range("D3").Value = $txt + '2.496.33'
$ws.Range("E3").Value = '  -0.72%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.06%  '

# Row 5: BNB
$ws.Range("D5").Value = $txt + '320.26'
$ws.Range("E5").Value = '  -0.88%  '

# Row 6: Solana
$ws.Range("D6").Value = $txt + '107.36'
$ws.Range("E6").Value = '  -2.09%  '

# Row 7: XRP
$ws.Range("E7").Value = '  -0.23%  '

# Row 8: USDC
$ws.Range("E8").Value = '  +0.03%  '

# Row 9: Cardano
$ws.Range("D9").Value = $txt + '0.537'
$ws.Range("E9").Value = '  -3.43%  '

# Row 10: Avalanche
$ws.Range("D10").Value = $txt + '39.41'
$ws.Range("E10").Value = '  -4.28%  '

# Row 11: Chainlink
$ws.Range("D11").Value = $txt + '20.20'
$ws.Range("E11").Value = '  +7.52%  '

# Row 12: Dogecoin
$ws.Range("D12").Value = $txt + '0.0812'
$ws.Range("E12").Value = '  -0.65%  '

# Row 13: TRON
$ws.Range("E13").Value = '  -0.46%  '

# Row 14: Polkadot
$ws.Range("D14").Value = $txt + '7.09'
$ws.Range("E14").Value = '  -2.70%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = $txt + '2.887.52'
$ws.Range("E15").Value = '  -0.43%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = $txt + '2.496.06'
$ws.Range("E16").Value = '  -0.66%  '

# Row 17: Polygon
$ws.Range("E17").Value = '  -2.47%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = $txt + '47.853.50'
$ws.Range("E18").Value = '  -0.05%  '

# Row 19: InternetComputer(DFINITY)
$ws.Range("D19").Value = $txt + '12.90'
$ws.Range("E19").Value = '  -3.50%  '

# Row 20: Uniswap
$ws.Range("D20").Value = $txt + '6.71'
$ws.Range("E20").Value = '  +0.73%  '

# Row 21: ShibaInu
$ws.Range("D21").Value = $txt + '0.0₃0938'
$ws.Range("E21").Value = '  -1.05%  '

# Row 22: ImmutableX
$ws.Range("D22").Value = $txt + '2.76'
$ws.Range("E22").Value = '  -2.47%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = $txt + '277.75'

# Row 24: Litecoin
$ws.Range("D24").Value = $txt + '71.46'
$ws.Range("E24").Value = '  +0.79%  '

# Row 25: PancakeSwap
$ws.Range("D25").Value = $txt + '2.53'
$ws.Range("E25").Value = '  -1.05%  '

# Row 26: Dai
$ws.Range("E26").Value = '  -0.11%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = $txt + '25.59'
$ws.Range("E27").Value = '  -1.60%  '

# Row 28: Toncoin (was Cosmos; list re-ranked)
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = $txt + '2.20'
$ws.Range("E28").Value = '  -0.03%  '

# Row 29: Cosmos (was Kaspa)
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = $txt + '9.72'
$ws.Range("E29").Value = '  -3.46%  '

# Row 30: Kaspa (was InjectiveProtocol)
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = $txt + '0.139'
$ws.Range("E30").Value = '  -0.37%  '

# Row 31: InjectiveProtocol (was Toncoin)
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = $txt + '34.93'
$ws.Range("E31").Value = '  -0.78%  '

# Row 32: OKB
$ws.Range("D32").Value = $txt + '49.40'
$ws.Range("E32").Value = '  -0.71%  '

# Row 33: Celestia
$ws.Range("D33").Value = $txt + '19.45'
$ws.Range("E33").Value = '  -3.49%  '

# Row 34: FirstDigitalUSD
$ws.Range("E34").Value = '  -0.10%  '

# Row 35: Filecoin
$ws.Range("D35").Value = $txt + '5.28'
$ws.Range("E35").Value = '  -2.07%  '

# Row 36: Hedera
$ws.Range("D36").Value = $txt + '0.0777'
$ws.Range("E36").Value = '  -1.52%  '

# Row 37: ARBITRUM
$ws.Range("E37").Value = '  -2.44%  '

# Row 38: RenderToken
$ws.Range("D38").Value = $txt + '4.60'
$ws.Range("E38").Value = '  -2.39%  '

# Row 39: LidoDAOToken
$ws.Range("E39").Value = '  -3.85%  '

# Row 40: Stellar
$ws.Range("E40").Value = '  -0.99%  '

# Row 41: Monero
$ws.Range("D41").Value = $txt + '120.68'
$ws.Range("E41").Value = '  +1.04%  '

# Row 42: WEMIXToken
$ws.Range("D42").Value = $txt + '2.21'
$ws.Range("E42").Value = '  -0.43%  '

# Row 43: EnergySwap
$ws.Range("D43").Value = $txt + '21.26'
$ws.Range("E43").Value = '  -6.21%  '

# Row 44: VeChain
$ws.Range("D44").Value = $txt + '0.0299'
$ws.Range("E44").Value = '  -0.10%  '

# Row 45: Maker
$ws.Range("D45").Value = $txt + '2.006.47'
$ws.Range("E45").Value = '  +0.14%  '

# Row 46: NEARProtocol
$ws.Range("D46").Value = $txt + '3.14'
$ws.Range("E46").Value = '  +1.76%  '

# Row 47: ApeXProtocol
$ws.Range("E47").Value = '  -1.92%  '

# Row 48: Stacks
$ws.Range("E48").Value = '  -0.67%  '

# Row 49: FraxShare
$ws.Range("D49").Value = $txt + '8.96'
$ws.Range("E49").Value = '  -1.42%  '

# Row 50: THORChain
$ws.Range("E50").Value = '  -1.52%  '

# Row 51: BitcoinSV
$ws.Range("D51").Value = $txt + '79.87'
$ws.Range("E51").Value = '  +1.97%  '
